$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "225 taka baki" -> "275 taka baki" (cell G24 holds this shared string)
$ws.Range("G24").Value = "275 taka baki"

# B3: 3188 -> 3268
$ws.Range("B3").Value = 3268

# F12: 0 -> 150
$ws.Range("F12").Value = 150

# F18: 75 -> 175
$ws.Range("F18").Value = 175

# Row 28 gets filled in with meal data for the day
$ws.Range("F28").Value = 160
$ws.Range("K28").Value = 2
$ws.Range("M28").Value = 2
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 2
$ws.Range("P28").Value = 2
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("T28").Value = 2

# Update view: move the active selection to D17 (25 tarikh / row for the 25th)
$ws.Activate()
$ws.Range("D17").Select()
